# #59: Check if all images exist, before importing
#
# The "Images" column (F) previously held a bogus placeholder value
# ("cintamani275_7") for every product row. After the image-existence
# check was added, each row's Images cell now holds the actual filename(s)
# found on disk: "001.JPG" for rows with a single image, and a
# comma-separated list ("001.JPG, 002.jpg") for the row with two images.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "001.JPG"
$ws.Range("F3").Value = "001.JPG"
$ws.Range("F4").Value = "001.JPG"
$ws.Range("F5").Value = "001.JPG"
$ws.Range("F6").Value = "001.JPG"
$ws.Range("F7").Value = "001.JPG, 002.jpg"

# Leave the selection where the author's last save left it: a single cell,
# F8 (just below the data table), rather than the old A4:XFD7 block select.
$ws.Range("F8").Select()
